# Generate Report for handoff
#
# The "86271fe5-5578-4006-ae23-7b5ea3b843f6" file has just been handed off
# (status "Ready for handoff" / reason "Ignored") in both the zh-cn and the
# de-de target-language sheets. Stamp the "Latest Handoff Datetime" cell
# (column D, row 6) on each of those sheets with the handoff timestamp.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D6").Value = "2016-01-18 03:38:08"
$wsDeDe.Range("D6").Value = "2016-01-18 03:38:19"
